$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new data row (GO338025225 / THAIS PEREIRA DE SOUSA) ---
# The new row must land right before the existing "GO338028636" row,
# which currently sits at worksheet row 51 (row 1 is the header).
$insertRow = 51
$ws.Rows.Item($insertRow).Insert()

$ws.Cells.Item($insertRow, 1).Value = "CENTRO OESTE"
$ws.Cells.Item($insertRow, 2).Value = "GOIANIA CENTRO NORTE"
$ws.Cells.Item($insertRow, 3).Value = "GO338025225"
$ws.Cells.Item($insertRow, 4).Value = "THAIS PEREIRA DE SOUSA"
$ws.Cells.Item($insertRow, 5).Value = 45127
$ws.Cells.Item($insertRow, 5).NumberFormat = "dd/MM/yyyy"
$ws.Cells.Item($insertRow, 6).Value = "VISA"
$ws.Cells.Item($insertRow, 7).Value = "62981695421"
# EMAIL column stays blank for this record; copy a blank cell in so the
# cell still exists (matches the empty <is><t/></is> cell in the source).
$ws.Cells.Item($insertRow + 1, 8).Copy($ws.Cells.Item($insertRow, 8))
$ws.Cells.Item($insertRow, 9).Value = "GABRIELA NUNES DA SILVA"
$ws.Cells.Item($insertRow, 10).Value = 6
$ws.Cells.Item($insertRow, 10).NumberFormat = "0"
$ws.Cells.Item($insertRow, 11).Value = 7
$ws.Cells.Item($insertRow, 11).NumberFormat = "0"
$ws.Cells.Item($insertRow, 12).Value = "CEP 74645190"

# --- 2) Remove the obsolete row for GO551002383 / MANOEL APARECIDO GONCALVES DE MOURA ---
$found = $ws.Cells.Find("GO551002383")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
